$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Oikeudet")
Write-Host $ws.Name
